$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Terminal Hortofrutícola Agro
# Chillán - Zanahoria". It belongs chronologically at the top of the
# existing price history block (row 126), so insert a new row there and
# let every subsequent record shift down by one (the former last record,
# previously on row 234, ends up on the new row 235).
$ws.Rows("126:126").Insert()

# Populate the newly inserted row with the new week's data.
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(126, 3).Value = "Ñuble"
$ws.Cells.Item(126, 4).Value = 44566
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = 100114013
$ws.Cells.Item(126, 7).Value = "Zanahoria"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 120
$ws.Cells.Item(126, 11).Value = 6500
$ws.Cells.Item(126, 12).Value = 7000
$ws.Cells.Item(126, 13).Value = 6750
$ws.Cells.Item(126, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(126, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(126, 16).Value = 338
$ws.Cells.Item(126, 17).Value = 20
$ws.Cells.Item(126, 18).Value = "Hortaliza"
